# Streamline presets, optional overrides, and killcam settings
# Replaces the "Balanced" preset label with "Standard" throughout the
# workbook, and updates a couple of explanatory note cells on the
# Overview sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview (Preset-First) sheet - unique note/header text updates
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview (Preset-First)")

$overview.Range("B3").Value = "Chance/Cooldown/Duration/Smoothness tables are derived from Intensity = Standard."
$overview.Range("B7").Value = "Killcam tables assume Third Person Distribution controls killcam."
$overview.Range("A69").Value = "Intensity Preset: Standard"
$overview.Range("A129").Value = "Chance Preset: Standard (Chance x1)"
$overview.Range("A169").Value = "Cooldown Preset: Standard (Cooldown x1)"

# ---------------------------------------------------------------------
# Per-trigger detail sheets - identical "Balanced" -> "Standard" pattern
# ---------------------------------------------------------------------
$triggerSheets = @("Basic Kill", "Critical", "Dismemberment", "Decapitation", "Parry", "Last Enemy", "Last Stand")

foreach ($sheetName in $triggerSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("B3").Value = "Chance/Cooldown/Duration/Smoothness derived from Intensity = Standard."
    $ws.Range("A11").Value = "Standard"
    $ws.Range("A16").Value = "Chance Presets (from Intensity Standard)"
    $ws.Range("A20").Value = "Standard"
    $ws.Range("A23").Value = "Cooldown Presets (from Intensity Standard)"
    $ws.Range("A27").Value = "Standard"
    $ws.Range("A31").Value = "Duration Presets (from Intensity Standard)"
    $ws.Range("A37").Value = "Smoothness Presets (from Intensity Standard)"
}
